$d = $word.ActiveDocument

$replacements = @(
    @{old="71÷2="; new="63÷8="},
    @{old="74÷4="; new="53÷7="},
    @{old="54÷9="; new="14÷8="},
    @{old="93÷5="; new="45÷6="},
    @{old="24÷6="; new="91÷3="},
    @{old="80÷7="; new="17÷7="},
    @{old="76÷6="; new="34÷2="},
    @{old="25÷6="; new="77÷9="},
    @{old="98÷3="; new="82÷4="},
    @{old="17÷4="; new="49÷8="},
    @{old="50÷3="; new="62÷2="},
    @{old="90÷5="; new="75÷4="},
    @{old="15÷2="; new="65÷6="},
    @{old="44÷8="; new="76÷9="},
    @{old="20÷4="; new="91÷3="},
    @{old="39÷5="; new="20÷5="},
    @{old="36÷3="; new="74÷7="},
    @{old="35÷2="; new="58÷3="},
    @{old="37÷4="; new="87÷9="},
    @{old="88÷8="; new="20÷8="},
    @{old="86÷2="; new="18÷2="},
    @{old="29÷8="; new="90÷6="},
    @{old="64÷2="; new="24÷7="},
    @{old="23÷4="; new="12÷3="},
    @{old="60÷4="; new="67÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
